$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update B2 total width value ---
$ws.Range("B2").Value = 325.50299999999999

# --- Row 3 (L01): update B3/C3 values ---
$ws.Range("B3").Value = 136.45400000000001
$ws.Range("C3").Value = 56.89

# --- Row 4 (L02): update B4/C4 values ---
$ws.Range("B4").Value = 136.45400000000001
$ws.Range("C4").Value = 59.79

# --- Row 5 (L03): change formatting to match rows 3/4 (fill/style 9/10), update values ---
$ws.Range("A4:E4").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)
$ws.Range("A5").Value = "L03"
$ws.Range("B5").Value = 136.45400000000001
$ws.Range("C5").Value = 70.319999999999993
$ws.Range("D5").Formula = "=D$2/B$2*B5"
$ws.Range("E5").Formula = "=E$2/C$2*C5"

# --- Row 6 (L04): update B6/C6 values ---
$ws.Range("B6").Value = 189.04900000000001
$ws.Range("C6").Value = 56.89

# --- Rows 7-8 (new L05/L06 rows): copy formatting from row 6, fill values/formulas ---
$ws.Range("A6:E6").Copy()
$ws.Range("A7:E7").PasteSpecial(-4122)
$ws.Range("A8:E8").PasteSpecial(-4122)

$ws.Range("A7").Value = "L05"
$ws.Range("B7").Value = 189.04900000000001
$ws.Range("C7").Value = 59.79
$ws.Range("D7:D8").Formula = "=D$2/B$2*B7"
$ws.Range("E7:E8").Formula = "=E$2/C$2*C7"

$ws.Range("A8").Value = "L06"
$ws.Range("B8").Value = 189.04900000000001
$ws.Range("C8").Value = 62.69
